$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Date"

# Copy C1's formatting (font) onto A1 and B1 so all header cells share style
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to match target state
$null = $ws.Range("D8").Select()
